# The commit "Fixed problems with reading size of classes and methods":
# on the methodNumberOfLines sheet, every method whose computed size was
# (incorrectly) reported as 0 lines should instead be reported as 1 line.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("methodNumberOfLines")

# Column C holds "Number of Lines" for each class/method row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 1) { $lastRow = 25 }

# Find every row whose "Number of Lines" cell currently reads "0".
$targetRows = @()
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Text -eq "0") {
        $targetRows += $r
    }
}

# Mark those cells as text first so the replacement "1" is written back as a
# string (matching how every other value in this column is stored as a
# shared string) instead of being auto-coerced into a numeric cell.
foreach ($r in $targetRows) {
    $ws.Cells.Item($r, 3).NumberFormat = "@"
}
foreach ($r in $targetRows) {
    $ws.Cells.Item($r, 3).Value = "1"
}
